$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
# Force text storage for numeric-looking price strings so Excel does not
# reinterpret them as numbers (which would drop meaningful trailing zeros
# and "thousands-dot" formatted values like "25.803.85").
$priceUpdates = @{
    'D2' = '25.803.85'
    'D3' = '1.630.60'
    'D4' = '0.997'
    'D5' = '214.09'
    'D8' = '0.255'
    'D10' = '19.66'
    'D11' = '0.0790'
    'D12' = '4.25'
    'D13' = '1.855.19'
    'D14' = '1.639.21'
    'D16' = '0.0₃0760'
    'D18' = '25.787.69'
    'D21' = '190.94'
    'D22' = '9.91'
    'D24' = '0.997'
    'D26' = '142.08'
    'D28' = '6.82'
    'D36' = '0.904'
    'D37' = '1.140.87'
    'D43' = '5.59'
    'D44' = '100.71'
    'D45' = '0.801'
    'D46' = '1.765.45'
    'D47' = '55.30'
    'D51' = '7.54'
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.ClearFormats()
}

# --- Volume(1h) (column E) updates ---
# These already contain non-numeric characters (spaces, "%"), so Excel
# keeps them as plain text without any extra coercion.
$volumeUpdates = @{
    'E2' = '  +0.22%  '
    'E3' = '  +0.18%  '
    'E4' = '  -0.54%  '
    'E5' = '  -0.19%  '
    'E6' = '  -0.09%  '
    'E7' = '  -0.54%  '
    'E8' = '  -0.86%  '
    'E9' = '  -0.68%  '
    'E10' = '  +0.38%  '
    'E11' = '  +0.89%  '
    'E12' = '  +0.18%  '
    'E13' = '  +0.17%  '
    'E14' = '  +0.10%  '
    'E15' = '  -0.09%  '
    'E16' = '  -0.39%  '
    'E17' = '  +0.06%  '
    'E18' = '  +0.14%  '
    'E19' = '  -0.52%  '
    'E20' = '  +0.11%  '
    'E21' = '  -1.60%  '
    'E22' = '  -0.04%  '
    'E23' = '  +0.78%  '
    'E24' = '  -0.63%  '
    'E25' = '  +1.82%  '
    'E26' = '  +1.83%  '
    'E27' = '  +1.78%  '
    'E28' = '  -0.32%  '
    'E29' = '  +0.06%  '
    'E30' = '  +0.21%  '
    'E31' = '  +1.47%  '
    'E32' = '  +0.08%  '
    'E33' = '  -0.67%  '
    'E34' = '  -0.17%  '
    'E35' = '  -0.11%  '
    'E36' = '  +0.85%  '
    'E37' = '  +2.84%  '
    'E38' = '  -0.04%  '
    'E39' = '  -2.16%  '
    'E40' = '  +0.15%  '
    'E41' = '  -0.67%  '
    'E42' = '  +0.15%  '
    'E43' = '  +0.72%  '
    'E44' = '  +0.75%  '
    'E45' = '  +0.47%  '
    'E46' = '  +0.43%  '
    'E47' = '  +0.66%  '
    'E48' = '  +7.20%  '
    'E49' = '  +2.20%  '
    'E50' = '  -0.27%  '
    'E51' = '  -1.41%  '
}

foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
